# Update 5 april 2020 — fill in the daily COVID-19 figures for each
# Sicilian province (columns C:H) that were still blank in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Data per province: Attualmente positivi, Ricoverati, Guariti, Morti,
# Casi totali, In isolamento domiciliare
$data = @{
    2  = @(104, 0,   2,  1,  107, 104)  # Agrigento
    3  = @(91,  23,  4,  7,  102, 68)   # Caltanissetta
    4  = @(525, 153, 23, 46, 594, 372)  # Catania
    5  = @(270, 170, 1,  13, 284, 100)  # Enna
    6  = @(314, 138, 15, 24, 353, 176)  # Messina
    7  = @(258, 73,  29, 12, 299, 185)  # Palermo
    8  = @(41,  7,   4,  3,  48,  34)   # Ragusa
    9  = @(77,  44,  25, 7,  109, 33)   # Siracusa
    10 = @(94,  24,  1,  3,  98,  70)   # Trapani
}

# Row 7 (Palermo) already had a lingering bold/header-styled blank cell in
# C7; clear that leftover formatting so it picks up the same plain style
# as the rest of the newly filled cells below.
$ws.Range("C7").ClearFormats()

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 3).Value = $values[0]
    $ws.Cells.Item($row, 4).Value = $values[1]
    $ws.Cells.Item($row, 5).Value = $values[2]
    $ws.Cells.Item($row, 6).Value = $values[3]
    $ws.Cells.Item($row, 7).Value = $values[4]
    $ws.Cells.Item($row, 8).Value = $values[5]
}

# Center the newly-filled block of figures.
$ws.Range("C2:H10").HorizontalAlignment = -4108

# Reflect the selection left on the sheet after the update.
$ws.Range("C2:H10").Select() | Out-Null
